$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H17").Value = 326392.97
$ws.Range("J17").Value = 326392.97
$ws.Range("L17").Value = 979178.9099999999
$ws.Range("N17").Value = -979514.9099999999
$ws.Range("H40").Value = 3691.75
$ws.Range("I40").Value = 2109.8333
$ws.Range("K40").Value = 2109.8333
$ws.Range("M40").Value = -1934.8333
$ws.Range("H48").Value = 1000
$ws.Range("I48").Value = 1000
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 3000
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -2708
$ws.Range("N48").ClearContents() | Out-Null
$ws.Range("H56").Value = 1000
$ws.Range("I56").Value = 1000
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 3000
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -2466
$ws.Range("N56").ClearContents() | Out-Null
$ws.Range("H62").Value = 4800
$ws.Range("I62").Value = 4720
$ws.Range("K62").Value = 4720
$ws.Range("M62").Value = -4096
$ws.Range("H65").Value = 4800
$ws.Range("I65").Value = 4720
$ws.Range("K65").Value = 23600
$ws.Range("M65").Value = -20480
$ws.Range("H125").Value = 3164.4614
$ws.Range("I125").Value = 1076.3334
$ws.Range("K125").Value = 9687.000599999999
$ws.Range("M125").Value = -7227.000599999999
$ws.Range("H132").Value = 29415326
$ws.Range("I132").Value = 34486676
$ws.Range("K132").Value = 103460028
$ws.Range("M132").Value = -103457498
$ws.Range("H138").Value = 106627.04
$ws.Range("I138").Value = 1254.6364
$ws.Range("J138").Value = 120104.91
$ws.Range("K138").Value = 3763.9092
$ws.Range("L138").Value = 360314.73
$ws.Range("M138").Value = 1376.0908
$ws.Range("N138").Value = -370594.73
$ws = $wb.Worksheets.Item(2)
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents() | Out-Null
$ws.Range("H5").Value = 3836.0715
$ws.Range("I5").Value = 4870
$ws.Range("J5").Value = 45
$ws.Range("K5").Value = 4870
$ws.Range("L5").Value = 45
$ws.Range("M5").Value = -4758
$ws.Range("N5").Value = -269
$ws.Range("H22").Value = 349.875
$ws.Range("I22").Value = 349.875
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 349.875
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -50.875
$ws.Range("N22").ClearContents() | Out-Null
$ws.Range("H32").Value = 6858.6562
$ws.Range("I32").Value = 6298.2393
$ws.Range("J32").Value = 19748.25
$ws.Range("K32").Value = 6298.2393
$ws.Range("L32").Value = 19748.25
$ws.Range("M32").Value = -6011.2393
$ws.Range("N32").Value = -20322.25
$ws.Range("H45").Value = 3356.2964
$ws.Range("I45").Value = 3374.8
$ws.Range("K45").Value = 3374.8
$ws.Range("M45").Value = -2997.8
$ws.Range("H61").Value = 11798.138
$ws.Range("I61").Value = 7890.346
$ws.Range("J61").Value = 45665.668
$ws.Range("K61").Value = 7890.346
$ws.Range("L61").Value = 45665.668
$ws.Range("M61").Value = -7678.346
$ws.Range("N61").Value = -46089.668
$ws.Range("H132").Value = 2430.9512
$ws.Range("I132").Value = 2218.0857
$ws.Range("K132").Value = 6654.257100000001
$ws.Range("M132").Value = -4124.257100000001
$ws.Range("H136").Value = 11798.138
$ws.Range("I136").Value = 7890.346
$ws.Range("J136").Value = 45665.668
$ws.Range("K136").Value = 23671.038
$ws.Range("L136").Value = 136997.004
$ws.Range("M136").Value = -21121.038
$ws.Range("N136").Value = -142097.004
$ws = $wb.Worksheets.Item(3)
$ws.Range("H4").Value = 3836.0715
$ws.Range("I4").Value = 4870
$ws.Range("J4").Value = 45
$ws.Range("K4").Value = 4870
$ws.Range("L4").Value = 45
$ws.Range("M4").Value = -4755
$ws.Range("N4").Value = -275
$ws = $wb.Worksheets.Item(4)
$ws.Range("H7").Value = 205
$ws.Range("I7").Value = 205
$ws.Range("K7").Value = 205
$ws.Range("M7").Value = -92
$ws.Range("H31").Value = 2697.946
$ws.Range("I31").Value = 2352.5862
$ws.Range("J31").Value = 3949.875
$ws.Range("K31").Value = 2352.5862
$ws.Range("L31").Value = 3949.875
$ws.Range("M31").Value = -2057.5862
$ws.Range("N31").Value = -4539.875
$ws.Range("H34").Value = 2697.946
$ws.Range("I34").Value = 2352.5862
$ws.Range("J34").Value = 3949.875
$ws.Range("K34").Value = 2352.5862
$ws.Range("L34").Value = 3949.875
$ws.Range("M34").Value = -2150.5862
$ws.Range("N34").Value = -4353.875
$ws.Range("H58").Value = 3277.7222
$ws.Range("I58").Value = 2886.4443
$ws.Range("J58").Value = 3669
$ws.Range("K58").Value = 2886.4443
$ws.Range("L58").Value = 3669
$ws.Range("M58").Value = -2683.4443
$ws.Range("N58").Value = -4075
$ws.Range("H68").Value = 50295
$ws.Range("J68").Value = 50295
$ws.Range("L68").Value = 50295
$ws.Range("N68").Value = -51793
$ws.Range("H71").Value = 50295
$ws.Range("J71").Value = 50295
$ws.Range("L71").Value = 150885
$ws.Range("N71").Value = -158373
$ws.Range("H132").Value = 13335709
$ws.Range("J132").Value = 1114
$ws.Range("L132").Value = 3342
$ws.Range("N132").Value = -8402
$ws.Range("H136").Value = 3277.7222
$ws.Range("I136").Value = 2886.4443
$ws.Range("J136").Value = 3669
$ws.Range("K136").Value = 8659.332900000001
$ws.Range("L136").Value = 11007
$ws.Range("M136").Value = -6109.332900000001
$ws.Range("N136").Value = -16107
$ws.Range("H141").Value = 636102.6
$ws.Range("J141").Value = 636102.6
$ws.Range("L141").Value = 636102.6
$ws.Range("N141").Value = -646462.6
$ws = $wb.Worksheets.Item(5)
$ws.Range("H2").Value = 117
$ws.Range("I2").Value = 50.5
$ws.Range("J2").Value = 250
$ws.Range("K2").Value = 303
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -190
$ws.Range("N2").Value = -1726
$ws.Range("H4").Value = 4778573.5
$ws.Range("I4").Value = 2556198.5
$ws.Range("J4").Value = 11128216
$ws.Range("K4").Value = 7668595.5
$ws.Range("L4").Value = 33384648
$ws.Range("M4").Value = -7668483.5
$ws.Range("N4").Value = -33384872
$ws.Range("H5").Value = 6598.9
$ws.Range("I5").Value = 940.6
$ws.Range("K5").Value = 2821.8
$ws.Range("M5").Value = -2709.8
$ws.Range("H132").Value = 3381.1428
$ws.Range("J132").Value = 3560.2258
$ws.Range("L132").Value = 32042.0322
$ws.Range("N132").Value = -37102.0322
$ws.Range("H135").Value = 6598.9
$ws.Range("I135").Value = 940.6
$ws.Range("K135").Value = 8465.4
$ws.Range("M135").Value = -5930.4
$ws = $wb.Worksheets.Item(6)
$ws.Range("H2").Value = 243.33333
$ws.Range("I2").Value = 188.57143
$ws.Range("K2").Value = 188.57143
$ws.Range("M2").Value = -75.57142999999999
$ws.Range("H122").Value = 3116.8462
$ws.Range("I122").Value = 1571.3334
$ws.Range("J122").Value = 6594.25
$ws.Range("K122").Value = 4714.0002
$ws.Range("L122").Value = 19782.75
$ws.Range("M122").Value = -2264.0002
$ws.Range("N122").Value = -24682.75
$ws.Range("H123").Value = 52999
$ws.Range("J123").Value = 52999
$ws.Range("L123").Value = 52999
$ws.Range("N123").Value = -57899
$ws.Range("H126").Value = 3753.0833
$ws.Range("I126").Value = 3226.3333
$ws.Range("K126").Value = 9678.999899999999
$ws.Range("M126").Value = -7208.999899999999
$ws.Range("H132").Value = 4307.125
$ws.Range("I132").Value = 4094.2666
$ws.Range("J132").Value = 7500
$ws.Range("K132").Value = 12282.7998
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -9752.799800000001
$ws.Range("N132").Value = -27560
$ws = $wb.Worksheets.Item(7)
$ws.Range("H46").Value = 10426.625
$ws.Range("I46").Value = 4798.8
$ws.Range("K46").Value = 4798.8
$ws.Range("M46").Value = -4610.8
$ws.Range("H55").Value = 380.7143
$ws.Range("I55").Value = 399.875
$ws.Range("K55").Value = 399.875
$ws.Range("M55").Value = -226.875
$ws.Range("H132").Value = 3219.625
$ws.Range("I132").Value = 3073.0386
$ws.Range("K132").Value = 9219.1158
$ws.Range("M132").Value = -6689.1158
$ws = $wb.Worksheets.Item(8)
$ws.Range("H126").Value = 2978.1
$ws.Range("I126").Value = 2099.6
$ws.Range("K126").Value = 6298.799999999999
$ws.Range("M126").Value = -3828.799999999999
$ws.Range("H132").Value = 2498.5
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents() | Out-Null
$ws.Range("H136").Value = 9349.759
$ws.Range("I136").Value = 15456.2
$ws.Range("J136").Value = 2807.1428
$ws.Range("K136").Value = 46368.60000000001
$ws.Range("L136").Value = 8421.428400000001
$ws.Range("M136").Value = -43818.60000000001
$ws.Range("N136").Value = -13521.4284
